$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4158.475
$ws.Range("I32").Value = 1966.5
$ws.Range("J32").Value = 4545.294
$ws.Range("K32").Value = 1966.5
$ws.Range("L32").Value = 4545.294
$ws.Range("M32").Value = -1640.5
$ws.Range("N32").Value = -5197.294

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 34506620
$ws.Range("I69").Value = 86871
$ws.Range("K69").Value = 260613
$ws.Range("M69").Value = -259739

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 34506620
$ws.Range("I72").Value = 86871
$ws.Range("K72").Value = 781839
$ws.Range("M72").Value = -777471

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 3298.6
$ws.Range("I82").Value = 3298.6
$ws.Range("K82").Value = 9895.799999999999
$ws.Range("M82").Value = -9489.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 3298.6
$ws.Range("I85").Value = 3298.6
$ws.Range("K85").Value = 9895.799999999999
$ws.Range("M85").Value = -8491.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1250.7273
$ws.Range("I92").Value = 1453.9259
$ws.Range("K92").Value = 1453.9259
$ws.Range("M92").Value = -205.9259

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1281.7391
$ws.Range("I135").Value = 944.7895
$ws.Range("K135").Value = 8503.1055
$ws.Range("M135").Value = -5968.1055

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1986.42
$ws.Range("I138").Value = 1035.8889
$ws.Range("K138").Value = 3107.6667
$ws.Range("M138").Value = 2032.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 11509.462
$ws.Range("I97").Value = 9083.237999999999
$ws.Range("K97").Value = 9083.237999999999
$ws.Range("M97").Value = -8587.237999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1040.2609
$ws.Range("I110").Value = 1109.0555
$ws.Range("K110").Value = 1109.0555
$ws.Range("M110").Value = 935.9445000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5719.1064
$ws.Range("I86").Value = 5597.615
$ws.Range("J86").Value = 5869.524
$ws.Range("K86").Value = 5597.615
$ws.Range("L86").Value = 5869.524
$ws.Range("M86").Value = -4474.615
$ws.Range("N86").Value = -8115.524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5719.1064
$ws.Range("I89").Value = 5597.615
$ws.Range("J89").Value = 5869.524
$ws.Range("K89").Value = 27988.075
$ws.Range("L89").Value = 29347.62
$ws.Range("M89").Value = -22372.075
$ws.Range("N89").Value = -40579.62

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1747.3
$ws.Range("I105").Value = 1611.4
$ws.Range("K105").Value = 1611.4
$ws.Range("M105").Value = 135.5999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 56827.8
$ws.Range("J130").Value = 56827.8
$ws.Range("L130").Value = 56827.8
$ws.Range("N130").Value = -66867.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 150000
$ws.Range("J131").Value = 150000
$ws.Range("L131").Value = 150000
$ws.Range("N131").Value = -160080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2947.175
$ws.Range("I134").Value = 964.3214
$ws.Range("K134").Value = 2892.9642
$ws.Range("M134").Value = -357.9642000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 198.3
$ws.Range("I7").Value = 135.375
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 135.375
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = -22.375
$ws.Range("N7").Value = -676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2612.3684
$ws.Range("I16").Value = 1945.3
$ws.Range("K16").Value = 1945.3
$ws.Range("M16").Value = -1658.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 9997.5
$ws.Range("J39").Value = 9997.5
$ws.Range("L39").Value = 9997.5
$ws.Range("N39").Value = -10779.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 9997.5
$ws.Range("J49").Value = 9997.5
$ws.Range("L49").Value = 9997.5
$ws.Range("N49").Value = -10361.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2329.7
$ws.Range("I62").Value = 1883.1666
$ws.Range("J62").Value = 2999.5
$ws.Range("K62").Value = 1883.1666
$ws.Range("L62").Value = 2999.5
$ws.Range("M62").Value = -1259.1666
$ws.Range("N62").Value = -4247.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2329.7
$ws.Range("I65").Value = 1883.1666
$ws.Range("J65").Value = 2999.5
$ws.Range("K65").Value = 9415.833000000001
$ws.Range("L65").Value = 14997.5
$ws.Range("M65").Value = -6295.833000000001
$ws.Range("N65").Value = -21237.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2612.3684
$ws.Range("I113").Value = 1945.3
$ws.Range("K113").Value = 1945.3
$ws.Range("M113").Value = 224.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2512
$ws.Range("I132").Value = 2231.6667
$ws.Range("K132").Value = 6695.000100000001
$ws.Range("M132").Value = -4165.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 25491.775
$ws.Range("I134").Value = 32342.666
$ws.Range("K134").Value = 97027.99800000001
$ws.Range("M134").Value = -94492.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1864.5
$ws.Range("I5").Value = 1241.1818
$ws.Range("K5").Value = 3723.5454
$ws.Range("M5").Value = -3611.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1233.9445
$ws.Range("I14").Value = 1233.9445
$ws.Range("K14").Value = 3701.8335
$ws.Range("M14").Value = -3528.8335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 131.44444
$ws.Range("J114").Value = 119.666664
$ws.Range("L114").Value = 358.999992
$ws.Range("N114").Value = -6866.999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4330
$ws.Range("I134").Value = 4330
$ws.Range("K134").Value = 12990
$ws.Range("M134").Value = -7920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1864.5
$ws.Range("I135").Value = 1241.1818
$ws.Range("K135").Value = 11170.6362
$ws.Range("M135").Value = -8635.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2373.7144
$ws.Range("I140").Value = 1769.5
$ws.Range("K140").Value = 5308.5
$ws.Range("M140").Value = -128.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3468.5
$ws.Range("I141").Value = 1877
$ws.Range("K141").Value = 5631
$ws.Range("M141").Value = -451

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2453.682
$ws.Range("I93").Value = 2493.7896
$ws.Range("J93").Value = 2199.6667
$ws.Range("K93").Value = 2493.7896
$ws.Range("L93").Value = 2199.6667
$ws.Range("M93").Value = -1245.7896
$ws.Range("N93").Value = -4695.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3202.2307
$ws.Range("I100").Value = 3007.2273
$ws.Range("J100").Value = 4274.75
$ws.Range("K100").Value = 3007.2273
$ws.Range("L100").Value = 4274.75
$ws.Range("M100").Value = -2466.2273
$ws.Range("N100").Value = -5356.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5132.123
$ws.Range("I132").Value = 4510.64
$ws.Range("J132").Value = 9571.286
$ws.Range("K132").Value = 13531.92
$ws.Range("L132").Value = 28713.858
$ws.Range("M132").Value = -11001.92
$ws.Range("N132").Value = -33773.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6619.836
$ws.Range("I62").Value = 2787.2
$ws.Range("J62").Value = 8489.415000000001
$ws.Range("K62").Value = 2787.2
$ws.Range("L62").Value = 8489.415000000001
$ws.Range("M62").Value = -2163.2
$ws.Range("N62").Value = -9737.415000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6619.836
$ws.Range("I65").Value = 2787.2
$ws.Range("J65").Value = 8489.415000000001
$ws.Range("K65").Value = 13936
$ws.Range("L65").Value = 42447.075
$ws.Range("M65").Value = -10816
$ws.Range("N65").Value = -48687.075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19979
$ws.Range("I107").Value = 13299
$ws.Range("K107").Value = 39897
$ws.Range("M107").Value = -37977

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1638.0312
$ws.Range("I126").Value = 1562.4615
$ws.Range("K126").Value = 4687.3845
$ws.Range("M126").Value = -2217.3845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8234.137000000001
$ws.Range("I132").Value = 9850.823
$ws.Range("J132").Value = 2737.4
$ws.Range("K132").Value = 29552.469
$ws.Range("L132").Value = 8212.200000000001
$ws.Range("M132").Value = -27022.469
$ws.Range("N132").Value = -13272.2

